# issue #5: add legislator_id, name, date into dataframe
# Adds three new columns (date, legislator_name, legislator_id) to the
# "股票" (stocks) worksheet, filling every existing data row with the
# filing date, the legislator's name, and their legislator id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$date = "2011-12-19"
$legislatorName = "柯建銘"
$legislatorId = 629

# ---- New header cells (row 1) ------------------------------------------
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Match the look of the existing header cells (bold, thin border, centered).
$headerRange = $ws.Range("H1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# ---- Fill every existing data row ---------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$hRange = $ws.Range("H2:H" + $lastRow)
$iRange = $ws.Range("I2:I" + $lastRow)
$jRange = $ws.Range("J2:J" + $lastRow)

# Force the date column to text first so "2011-12-19" isn't reinterpreted
# as a date serial number, then fill every data row in one shot.
$hRange.NumberFormat = "@"
$hRange.Value = $date
$iRange.Value = $legislatorName
$jRange.Value = $legislatorId
